# "prueba con dos barridos" - append a second sweep of sensor readings
# (rows 177-352) below the existing data (rows 1-176), mirroring the
# angle column (A) and appending the new distance readings (B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aVals = @(5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,71,72,73,74,75,76,77,78,79,80,81,82,83,84,85,86,87,88,89,90,91,92,93,94,95,96,97,98,99,100,101,102,103,104,105,106,107,108,109,110,111,112,113,114,115,116,117,118,119,120,121,122,123,124,125,126,127,128,129,130,131,132,133,134,135,136,137,138,139,140,141,142,143,144,145,146,147,148,149,150,151,152,153,154,155,156,157,158,159,160,161,162,163,164,165,166,167,168,169,170,171,172,173,174,175,176,177,178,179,180)
$bVals = @(80,16,16,16,16,16,16,16,17,16,80,16,80,16,16,16,16,16,16,16,16,16,16,16,16,16,17,16,16,16,16,16,16,17,17,17,17,17,17,17,18,16,17,16,18,17,18,17,19,18,28,28,27,80,27,27,28,28,28,28,28,29,26,30,30,31,80,80,80,80,80,21,23,21,22,21,22,23,22,22,22,21,21,21,21,21,21,21,21,21,21,21,21,21,21,21,21,21,20,21,20,21,20,21,21,21,21,21,21,21,21,21,22,22,22,22,23,23,26,25,26,25,26,23,25,25,25,25,25,25,25,25,25,25,25,25,26,25,25,26,25,26,26,26,26,26,27,27,28,28,29,34,80,80,15,80,14,14,14,14,13,13,13,12,13,12,12,12,12,12,12,12,12,1,12,12)

$startRow = 177
$endRow = $startRow + $aVals.Length - 1

$data = New-Object 'object[,]' $aVals.Length,2
for ($i = 0; $i -lt $aVals.Length; $i++) {
    $data[$i,0] = $aVals[$i]
    $data[$i,1] = $bVals[$i]
}

$targetRange = $ws.Range($ws.Cells.Item($startRow,1), $ws.Cells.Item($endRow,2))
$targetRange.Value = $data

$ws.Application.ActiveWindow.ScrollRow = 330
$ws.Range("A308:B352").Select()
